$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 402.64706
$ws.Range("I53").Value = 482
$ws.Range("J53").Value = 212.2
$ws.Range("K53").Value = 482
$ws.Range("L53").Value = 212.2
$ws.Range("M53").Value = 155
$ws.Range("N53").Value = -1486.2

$ws.Range("H54").Value = 13353.083
$ws.Range("I54").Value = 6500
$ws.Range("J54").Value = 14723.7
$ws.Range("K54").Value = 6500
$ws.Range("L54").Value = 14723.7
$ws.Range("M54").Value = -6014
$ws.Range("N54").Value = -15695.7

$ws.Range("H55").Value = 378.42856
$ws.Range("I55").Value = 285
$ws.Range("J55").Value = 415.8
$ws.Range("K55").Value = 285
$ws.Range("L55").Value = 415.8
$ws.Range("M55").Value = -71
$ws.Range("N55").Value = -843.8

$ws.Range("H112").Value = 1038
$ws.Range("I112").Value = 666.6667
$ws.Range("J112").Value = 1093.7
$ws.Range("K112").Value = 2000.0001
$ws.Range("L112").Value = 3281.1
$ws.Range("M112").Value = -892.0001
$ws.Range("N112").Value = -5497.1

$ws.Range("H113").Value = 2638.8462
$ws.Range("I113").Value = 1568.3334
$ws.Range("J113").Value = 2960
$ws.Range("K113").Value = 1568.3334
$ws.Range("L113").Value = 2960
$ws.Range("M113").Value = 1685.6666
$ws.Range("N113").Value = -9468

$ws.Range("H115").Value = 480.625
$ws.Range("I115").Value = 335
$ws.Range("J115").Value = 1500
$ws.Range("K115").Value = 1005
$ws.Range("L115").Value = 4500
$ws.Range("M115").Value = 562
$ws.Range("N115").Value = -7634

$ws.Range("H127").Value = 843.2857
$ws.Range("I127").Value = 611.36365
$ws.Range("J127").Value = 1098.4
$ws.Range("K127").Value = 1834.09095
$ws.Range("L127").Value = 3295.2
$ws.Range("M127").Value = 3125.90905
$ws.Range("N127").Value = -13215.2

$ws.Range("H129").Value = 6412.974
$ws.Range("J129").Value = 8473.311
$ws.Range("L129").Value = 25419.933
$ws.Range("N129").Value = -35419.933

$ws.Range("H138").Value = 2438.7183
$ws.Range("J138").Value = 3127.6978
$ws.Range("L138").Value = 9383.0934
$ws.Range("N138").Value = -19663.0934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 994.3125
$ws.Range("I2").Value = 933.8333
$ws.Range("J2").Value = 1030.6
$ws.Range("K2").Value = 933.8333
$ws.Range("L2").Value = 1030.6
$ws.Range("M2").Value = -820.8333
$ws.Range("N2").Value = -1256.6

$ws.Range("H7").Value = 31000
$ws.Range("J7").Value = 31000
$ws.Range("L7").Value = 31000
$ws.Range("N7").Value = -31228

$ws.Range("H45").Value = 2346.5454
$ws.Range("I45").Value = 4404
$ws.Range("J45").Value = 1575
$ws.Range("K45").Value = 4404
$ws.Range("L45").Value = 1575
$ws.Range("M45").Value = -4027
$ws.Range("N45").Value = -2329

$ws.Range("H59").Value = 40059
$ws.Range("J59").Value = 40059
$ws.Range("L59").Value = 40059
$ws.Range("N59").Value = -41667

$ws.Range("H63").Value = 2715.75
$ws.Range("I63").Value = 2268.9
$ws.Range("J63").Value = 4950
$ws.Range("K63").Value = 2268.9
$ws.Range("L63").Value = 4950
$ws.Range("M63").Value = -1582.9
$ws.Range("N63").Value = -6322

$ws.Range("H66").Value = 2715.75
$ws.Range("I66").Value = 2268.9
$ws.Range("J66").Value = 4950
$ws.Range("K66").Value = 11344.5
$ws.Range("L66").Value = 24750
$ws.Range("M66").Value = -7912.5
$ws.Range("N66").Value = -31614

$ws.Range("H110").Value = 3471.8572
$ws.Range("I110").Value = 2327.75
$ws.Range("J110").Value = 4997.3335
$ws.Range("K110").Value = 2327.75
$ws.Range("L110").Value = 4997.3335
$ws.Range("M110").Value = -282.75
$ws.Range("N110").Value = -9087.333500000001

$ws.Range("H116").Value = 994.3125
$ws.Range("I116").Value = 933.8333
$ws.Range("J116").Value = 1030.6
$ws.Range("K116").Value = 933.8333
$ws.Range("L116").Value = 1030.6
$ws.Range("M116").Value = 1360.1667
$ws.Range("N116").Value = -5618.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 994.3125
$ws.Range("I3").Value = 933.8333
$ws.Range("J3").Value = 1030.6
$ws.Range("K3").Value = 933.8333
$ws.Range("L3").Value = 1030.6
$ws.Range("M3").Value = -819.8333
$ws.Range("N3").Value = -1258.6

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H82").Value = 24195.182
$ws.Range("I82").Value = 14752.333
$ws.Range("J82").Value = 27736.25
$ws.Range("K82").Value = 14752.333
$ws.Range("L82").Value = 27736.25
$ws.Range("M82").Value = -14369.333
$ws.Range("N82").Value = -28502.25

$ws.Range("H85").Value = 24195.182
$ws.Range("I85").Value = 14752.333
$ws.Range("J85").Value = 27736.25
$ws.Range("K85").Value = 14752.333
$ws.Range("L85").Value = 27736.25
$ws.Range("M85").Value = -13426.333
$ws.Range("N85").Value = -30388.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 19800
$ws.Range("J59").Value = 19800
$ws.Range("L59").Value = 19800
$ws.Range("N59").Value = -22090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132697.33
$ws.Range("I2").Value = 220016.89
$ws.Range("J2").Value = 1718
$ws.Range("K2").Value = 1320101.34
$ws.Range("L2").Value = 10308
$ws.Range("M2").Value = -1319988.34
$ws.Range("N2").Value = -10534

$ws.Range("H50").Value = 263.85715
$ws.Range("I50").Value = 35.5
$ws.Range("J50").Value = 568.3333
$ws.Range("K50").Value = 106.5
$ws.Range("L50").Value = 1704.9999
$ws.Range("M50").Value = 374.5
$ws.Range("N50").Value = -2666.9999

$ws.Range("H53").Value = 263.85715
$ws.Range("I53").Value = 35.5
$ws.Range("J53").Value = 568.3333
$ws.Range("K53").Value = 106.5
$ws.Range("L53").Value = 1704.9999
$ws.Range("M53").Value = 374.5
$ws.Range("N53").Value = -2666.9999

$ws.Range("H54").Value = 3829.8333
$ws.Range("J54").Value = 3829.8333
$ws.Range("L54").Value = 11489.4999
$ws.Range("N54").Value = -12607.4999

$ws.Range("H55").Value = 1686.8125
$ws.Range("J55").Value = 1686.8125
$ws.Range("L55").Value = 5060.4375
$ws.Range("N55").Value = -5414.4375

$ws.Range("H68").Value = 766.36
$ws.Range("I68").Value = 554.61536
$ws.Range("J68").Value = 995.75
$ws.Range("K68").Value = 1663.84608
$ws.Range("L68").Value = 2987.25
$ws.Range("M68").Value = -852.84608
$ws.Range("N68").Value = -4609.25

$ws.Range("H71").Value = 766.36
$ws.Range("I71").Value = 554.61536
$ws.Range("J71").Value = 995.75
$ws.Range("K71").Value = 4991.53824
$ws.Range("L71").Value = 8961.75
$ws.Range("M71").Value = -935.5382399999999
$ws.Range("N71").Value = -17073.75

$ws.Range("H131").Value = 913.5484
$ws.Range("J131").Value = 1061.6
$ws.Range("L131").Value = 3184.8
$ws.Range("N131").Value = -13264.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 10000
$ws.Range("J53").Value = 10000
$ws.Range("L53").Value = 10000
$ws.Range("N53").Value = -11262

$ws.Range("H55").Value = 2666.3333
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 2999.5
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 2999.5
$ws.Range("M55").Value = -1673
$ws.Range("N55").Value = -3653.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 839.1070999999999
$ws.Range("I46").Value = 727.36365
$ws.Range("J46").Value = 911.41174
$ws.Range("K46").Value = 727.36365
$ws.Range("L46").Value = 911.41174
$ws.Range("M46").Value = -539.36365
$ws.Range("N46").Value = -1287.41174

$ws.Range("H55").Value = 269.73914
$ws.Range("J55").Value = 235.33333
$ws.Range("L55").Value = 235.33333
$ws.Range("N55").Value = -581.3333299999999

Write-Output "edit applied"
